$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 346.15384
$ws.Range("I12").Value = 270
$ws.Range("J12").Value = 380
$ws.Range("K12").Value = 270
$ws.Range("L12").Value = 380
$ws.Range("M12").Value = -100
$ws.Range("N12").Value = -720
$ws.Range("H32").Value = 4546.875
$ws.Range("I32").Value = 4125
$ws.Range("J32").Value = 4687.5
$ws.Range("K32").Value = 4125
$ws.Range("L32").Value = 4687.5
$ws.Range("M32").Value = -3799
$ws.Range("N32").Value = -5339.5
$ws.Range("H64").Value = 4531.25
$ws.Range("I64").Value = 3750
$ws.Range("K64").Value = 3750
$ws.Range("M64").Value = -3502
$ws.Range("H67").Value = 4531.25
$ws.Range("I67").Value = 3750
$ws.Range("K67").Value = 3750
$ws.Range("M67").Value = -2892
$ws.Range("H116").Value = 4848.467
$ws.Range("I116").Value = 4679.778
$ws.Range("J116").Value = 5101.5
$ws.Range("K116").Value = 4679.778
$ws.Range("L116").Value = 5101.5
$ws.Range("M116").Value = -1237.778
$ws.Range("N116").Value = -11985.5
$ws.Range("H132").Value = 953.4286
$ws.Range("I132").Value = 834.44116
$ws.Range("K132").Value = 2503.32348
$ws.Range("M132").Value = 26.67651999999998
$ws.Range("H135").Value = 659.5
$ws.Range("I135").Value = 659.5
$ws.Range("K135").Value = 5935.5
$ws.Range("M135").Value = -3400.5
$ws.Range("H137").Value = 7463.154
$ws.Range("J137").Value = 14998.333
$ws.Range("L137").Value = 44994.999
$ws.Range("N137").Value = -50094.999
$ws.Range("H138").Value = 2513.8767
$ws.Range("J138").Value = 2791.5173
$ws.Range("L138").Value = 8374.5519
$ws.Range("N138").Value = -18654.5519
$ws.Range("H141").Value = 10593.211
$ws.Range("I141").Value = 9690
$ws.Range("K141").Value = 29070
$ws.Range("M141").Value = -23890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 33335434
$ws.Range("I45").Value = 38463316
$ws.Range("K45").Value = 38463316
$ws.Range("M45").Value = -38462939
$ws.Range("H74").Value = 6501163.5
$ws.Range("I74").Value = 7143993.5
$ws.Range("K74").Value = 7143993.5
$ws.Range("M74").Value = -7143119.5
$ws.Range("H77").Value = 6501163.5
$ws.Range("I77").Value = 7143993.5
$ws.Range("K77").Value = 35719967.5
$ws.Range("M77").Value = -35715599.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 648.25
$ws.Range("I5").Value = 764.3333
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 764.3333
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -651.3333
$ws.Range("N5").Value = -526
$ws.Range("H99").Value = 2908.8948
$ws.Range("I99").Value = 2157.9092
$ws.Range("J99").Value = 3941.5
$ws.Range("K99").Value = 2157.9092
$ws.Range("L99").Value = 3941.5
$ws.Range("M99").Value = -659.9092000000001
$ws.Range("N99").Value = -6937.5
$ws.Range("H105").Value = 2612.6956
$ws.Range("I105").Value = 1515
$ws.Range("J105").Value = 2843.7896
$ws.Range("K105").Value = 1515
$ws.Range("L105").Value = 2843.7896
$ws.Range("M105").Value = 232
$ws.Range("N105").Value = -6337.7896
$ws.Range("H134").Value = 526074.6
$ws.Range("I134").Value = 1864.0555
$ws.Range("K134").Value = 5592.166499999999
$ws.Range("M134").Value = -3057.166499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1672024.6
$ws.Range("H65").Value = 1672024.6
$ws.Range("H132").Value = 6621.615
$ws.Range("I132").Value = 4235.278
$ws.Range("K132").Value = 12705.834
$ws.Range("M132").Value = -10175.834
$ws.Range("H134").Value = 4037.818
$ws.Range("I134").Value = 3379.5557
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 10138.6671
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -7603.667099999999
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 11212.583
$ws.Range("J134").Value = 12779.613
$ws.Range("L134").Value = 38338.839
$ws.Range("N134").Value = -48478.839
$ws.Range("H136").Value = 11354.5
$ws.Range("I136").Value = 5709.25
$ws.Range("J136").Value = 16999.75
$ws.Range("K136").Value = 17127.75
$ws.Range("L136").Value = 50999.25
$ws.Range("M136").Value = -12027.75
$ws.Range("N136").Value = -61199.25
$ws.Range("H137").Value = 5257.9443
$ws.Range("I137").Value = 4600.5
$ws.Range("J137").Value = 5586.6665
$ws.Range("K137").Value = 13801.5
$ws.Range("L137").Value = 16759.9995
$ws.Range("M137").Value = -8701.5
$ws.Range("N137").Value = -26959.9995
$ws.Range("H138").Value = 4502.727
$ws.Range("H139").Value = 4110.6
$ws.Range("I139").Value = 3704.375
$ws.Range("K139").Value = 11113.125
$ws.Range("M139").Value = -5973.125
$ws.Range("H140").Value = 234186.77
$ws.Range("I140").Value = 234186.77
$ws.Range("K140").Value = 702560.3099999999
$ws.Range("M140").Value = -697380.3099999999
$ws.Range("H141").Value = 176300
$ws.Range("J141").Value = 14866.667
$ws.Range("L141").Value = 44600.001
$ws.Range("N141").Value = -54960.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24332.166
$ws.Range("I70").Value = 24332.166
$ws.Range("K70").Value = 24332.166
$ws.Range("M70").Value = -24062.166
$ws.Range("H73").Value = 24332.166
$ws.Range("I73").Value = 24332.166
$ws.Range("K73").Value = 24332.166
$ws.Range("M73").Value = -23396.166
$ws.Range("H97").Value = 1426
$ws.Range("I97").Value = 1322.1818
$ws.Range("J97").Value = 1997
$ws.Range("K97").Value = 1322.1818
$ws.Range("L97").Value = 1997
$ws.Range("M97").Value = -826.1818000000001
$ws.Range("N97").Value = -2989
$ws.Range("H102").Value = 2405.6316
$ws.Range("I102").Value = 1866.1852
$ws.Range("K102").Value = 1866.1852
$ws.Range("M102").Value = -244.1851999999999
$ws.Range("H126").Value = 2836
$ws.Range("I126").Value = 2015.0625
$ws.Range("J126").Value = 4712.4287
$ws.Range("K126").Value = 6045.1875
$ws.Range("L126").Value = 14137.2861
$ws.Range("M126").Value = -3575.1875
$ws.Range("N126").Value = -19077.2861
$ws.Range("H132").Value = 90919760
$ws.Range("I132").Value = 100001736
$ws.Range("K132").Value = 300005208
$ws.Range("M132").Value = -300002678
$ws.Range("H136").Value = 10473.25
$ws.Range("J136").Value = 10473.25
$ws.Range("L136").Value = 31419.75
$ws.Range("N136").Value = -36519.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 29413016
$ws.Range("I93").Value = 47619990
$ws.Range("J93").Value = 1742.6923
$ws.Range("K93").Value = 47619990
$ws.Range("L93").Value = 1742.6923
$ws.Range("M93").Value = -47618742
$ws.Range("N93").Value = -4238.6923
$ws.Range("H100").Value = 4216.6665
$ws.Range("J100").Value = 4500
$ws.Range("L100").Value = 4500
$ws.Range("N100").Value = -5582
$ws.Range("H132").Value = 1002195.44
$ws.Range("I132").Value = 126018.125
$ws.Range("K132").Value = 378054.375
$ws.Range("M132").Value = -375524.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 12248.375
$ws.Range("I61").Value = 7641.7144
$ws.Range("K61").Value = 7641.7144
$ws.Range("M61").Value = -7349.7144
$ws.Range("H107").Value = 35715056
$ws.Range("I107").Value = 41667348
$ws.Range("K107").Value = 125002044
$ws.Range("M107").Value = -125000124
